$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.152.02'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '2.210.56'
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '295.21'
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.63'
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -0.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.85'
$ws.Range("E10").Value = '  +1.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '51.40'
$ws.Range("E11").Value = '  +5.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0780'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  +2.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.38'
$ws.Range("E14").Value = '  -1.67%  '
$ws.Range("D15").Value = '2.555.70'
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.80'
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("D17").Value = '2.254.14'
$ws.Range("E17").Value = '  +0.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.736'
$ws.Range("E18").Value = '  +1.17%  '
$ws.Range("D19").Value = '40.065.41'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.27'
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("E22").Value = '  -1.31%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.73'
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("E26").Value = '  +1.06%  '
$ws.Range("E27").Value = '  -1.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.15'
$ws.Range("E28").Value = '  +1.95%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.32'
$ws.Range("E29").Value = '  +1.26%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.08'
$ws.Range("E30").Value = '  -4.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.40'
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.07'
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.95'
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.02'
$ws.Range("E35").Value = '  +4.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0714'
$ws.Range("E36").Value = '  -0.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.32'
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("E38").Value = '  +1.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.101'
$ws.Range("E39").Value = '  +2.79%  '
$ws.Range("E40").Value = '  +2.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.67'
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").Value = '2.074.12'
$ws.Range("E42").Value = '  -2.06%  '
$ws.Range("E43").Value = '  -1.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.43'
$ws.Range("E44").Value = '  +7.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0271'
$ws.Range("E45").Value = '  +1.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.98'
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.81'
$ws.Range("E47").Value = '  +5.47%  '
$ws.Range("E48").Value = '  -10.64%  '
$ws.Range("D49").Value = '2.427.50'
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.12'
$ws.Range("E50").Value = '  +1.82%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.47'
$ws.Range("E51").Value = '  +0.81%  '
